$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.886.78"
$ws.Range("E2").Value = "  -4.40%  "
$ws.Range("D3").Value = "3.226.80"
$ws.Range("E3").Value = "  -5.53%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "173.76"
$ws.Range("E5").Value = "  -5.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "515.76"
$ws.Range("E6").Value = "  -3.71%  "
$ws.Range("E7").Value = "  -4.27%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "3.228.50"
$ws.Range("E9").Value = "  -5.31%  "
$ws.Range("E10").Value = "  -5.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.43"
$ws.Range("E11").Value = "  -10.23%  "
$ws.Range("E12").Value = "  -4.16%  "
$ws.Range("E13").Value = "  -2.66%  "
$ws.Range("E14").Value = "  -5.86%  "
$ws.Range("D15").Value = "3.746.49"
$ws.Range("E15").Value = "  -5.36%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.225.37"
$ws.Range("E16").Value = "  -5.43%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.115"
$ws.Range("E17").Value = "  -6.76%  "
$ws.Range("D18").Value = "62.885.71"
$ws.Range("E18").Value = "  -4.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.09"
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.96"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("E21").Value = "  -3.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "365.14"
$ws.Range("E22").Value = "  -4.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.71"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.16"
$ws.Range("E24").Value = "  -4.15%  "
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  +5.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.10"
$ws.Range("E27").Value = "  +3.98%  "
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.18"
$ws.Range("E29").Value = "  -4.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.12"
$ws.Range("E30").Value = "  -5.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "654.71"
$ws.Range("E31").Value = "  -5.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.14"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.29"
$ws.Range("E33").Value = "  -7.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.05"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.32"
$ws.Range("E36").Value = "  -7.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.38"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("E39").Value = "  -5.02%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0717"
$ws.Range("E40").Value = "  +14.29%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("D43").Value = "2.855.40"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.49"
$ws.Range("E44").Value = "  +5.17%  "
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0388"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.58"
$ws.Range("E47").Value = "  -7.65%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.81"
$ws.Range("E48").Value = "  +8.57%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.95"
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.73"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.123"
$ws.Range("E51").Value = "  -3.29%  "
